# Issue #29 add magnetometer to parts list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing hyperlinks (ref cell + target URL) in their current
# order *before* we touch the sheet, so we can recreate them afterwards
# pointing at the row that the insert will shift them to.
$oldLinks = @()
for ($i = 1; $i -le $ws.Hyperlinks.Count; $i++) {
    $hl = $ws.Hyperlinks.Item($i)
    $oldLinks += , @($hl.Range.Address($false, $false), $hl.Address)
}

# Remove all hyperlinks up front -- Rows.Insert() does not renumber the
# <hyperlinks> refs on its own, so we rebuild them from scratch below.
$ws.Hyperlinks.Delete()

# Insert a brand-new row above the old row 12 (GPS), pushing GPS and
# everything below it down by one row.
$ws.Rows.Item(12).Insert()

# The inserted row copies formatting from the row above (the "Sensors"
# section header); reset it back to plain/default formatting to match
# the rest of the data rows.
$ws.Range("A12").Style = "Normal"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Style = "Hyperlink"

# Fill in the new Magnetometer row. Write column C before column A so the
# shared-string table picks up the URL before the label (matches source order).
$ws.Range("C12").Value2 = "https://www.sparkfun.com/products/10530"
$ws.Range("A12").Value2 = "Magnetometer"
$ws.Range("B12").Value2 = 1

# Re-create the hyperlinks. Any ref at row >= 12 needs to shift down by one
# to stay attached to the same logical row after the insert.
foreach ($link in $oldLinks) {
    $ref = $link[0]
    $target = $link[1]

    if ($ref -match '^([A-Z]+)(\d+)$') {
        $col = $Matches[1]
        $row = [int]$Matches[2]
        if ($row -ge 12) {
            $row = $row + 1
        }
        $newRef = "$col$row"
    } else {
        $newRef = $ref
    }

    $ws.Hyperlinks.Add($ws.Range($newRef), $target) | Out-Null
}

# Finally add the hyperlink for the new Magnetometer row itself.
$ws.Hyperlinks.Add($ws.Range("C12"), "https://www.sparkfun.com/products/10530") | Out-Null

$ws.Range("A9").Select()
